# Update countries & provincias Spain
# Applies the data refresh described in the commit:
#  - Update the "last updated" timestamp
#  - Refresh case statistics for several countries (rows reference the
#    sheet's sorted-by-total-cases order, which causes Sudafrica and
#    Egipto to swap places once Sudafrica's totals overtake Egipto's)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Header timestamp
$ws.Range("A1").Value = "Datos actualizados a 30 de Abril de 2020 a las 20:52"

# Estados Unidos (row 4)
$ws.Range("E4").Value = 866247
$ws.Range("G4").Value = 888
$ws.Range("H4").Value = 62543

# Francia (row 8)
$ws.Range("D8").Value = 49476
$ws.Range("E8").Value = 92568
$ws.Range("F8").Value = 4019
$ws.Range("G8").Value = 289
$ws.Range("H8").Value = 24376

# Sudafrica now overtakes Egipto in total cases, so they swap rows (52/53)
$ws.Range("A52").Value = "Sudafrica"
$ws.Range("B52").Value = 5647
$ws.Range("C52").Value = 297
$ws.Range("D52").Value = 2073
$ws.Range("E52").Value = 3471
$ws.Range("F52").Value = 36
$ws.Range("G52").Value = 0
$ws.Range("H52").Value = 103

$ws.Range("A53").Value = "Egipto"
$ws.Range("B53").Value = 5537
$ws.Range("C53").Value = 269
$ws.Range("D53").Value = 1381
$ws.Range("E53").Value = 3764
$ws.Range("F53").Value = 0
$ws.Range("G53").Value = 12
$ws.Range("H53").Value = 392

# Uzbekistan (row 72)
$ws.Range("B72").Value = 2039
$ws.Range("C72").Value = 37
$ws.Range("E72").Value = 897

# Sri Lanka (row 103)
$ws.Range("B103").Value = 663
$ws.Range("C103").Value = 14
$ws.Range("D103").Value = 154
$ws.Range("E103").Value = 502

# Mali (row 112)
$ws.Range("B112").Value = 490
$ws.Range("C112").Value = 8
$ws.Range("D112").Value = 135
$ws.Range("E112").Value = 329
$ws.Range("G112").Value = 1
$ws.Range("H112").Value = 26

# Libia (row 165)
$ws.Range("E165").Value = 40
$ws.Range("G165").Value = 1
$ws.Range("H165").Value = 3
